$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.259.38"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.176.52"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.77"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.58"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.170.77"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  -6.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.66"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "3.703.03"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "66.356.22"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.36"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "3.184.68"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.56"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.39"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.08"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.83"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.31"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  +6.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("E30").Value = "  +7.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.99"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.52"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "509.40"
$ws.Range("E36").Value = "  +5.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.71"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0892"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("E40").Value = "  +6.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.81"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.299"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").Value = "0.0₃0667"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "2.841.88"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.23"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.30"
$ws.Range("E51").Value = "  +5.04%  "
